# Applies the "Freshly forged Ammo" augment + "Dashing" synergy addition
# described in the commit message / diff.

$wb = $excel.ActiveWorkbook
$wsAug = $wb.Worksheets.Item("Augments")
$wsSyn = $wb.Worksheets.Item("Synergies")

# ----------------------------------------------------------------------
# Sheet "Augments" (sheet1)
# ----------------------------------------------------------------------

# Row 3 ("Bouncy Balls") gains a wrapped Code cell and a taller row.
$wsAug.Range("G3").WrapText = $true
$wsAug.Rows.Item(3).RowHeight = 60

# Row 13 ("Personal Space") synergies list now also references the new
# "Dashing" synergy (id 5) on the Synergies sheet.
$wsAug.Range("D13").Value = "3,1, 5"

# New row 14: "Freshly Forged Ammo" augment.
$wsAug.Range("A14").Value = 12
$wsAug.Range("B14").Value = "Freshly Forged Ammo"
$wsAug.Range("C14").Value = 3
$wsAug.Range("D14").Value = 0
$wsAug.Range("E14").Value = "icons/FreshAmmo.png"

$wsAug.Range("F14").Value = "<b>+<color=#c5c5c5ff><color=#00d100>50%</color>/70%/100%</color></b> Damage increase`n<b>+<color=#c5c5c5ff><color=#00d100>+1 Bullet, +1 Grenade</color>/+1 Bullet, +2 Grenade/+2 Bullet, +2 Grenade</color></b> `nThe first shot after a reload does way more damage."
$wsAug.Range("G14").Value = "def OnAttached() {   `nAddModifier(`"Character`", `"extraDamageAfterReload`", `"Flat`", 0.5);`nAddModifier(`"Bullet`", `"maxClip`", `"Flat`", 1);`nAddModifier(`"Grenade`", `"maxClip`", `"Flat`",1);`n} "
$wsAug.Range("H14").Value = "<b>+<color=#c5c5c5ff>50%/<color=#00d100>70%</color>/100%</color></b> Damage increase`n<b>+<color=#c5c5c5ff>+1 Bullet, +1 Grenade/<color=#00d100>+1 Bullet, +2 Grenade</color>/+2 Bullet, +2 Grenade</color></b> `nThe first shot after a reload does way more damage."
$wsAug.Range("I14").Value = "def OnAttached() {   `nAddModifier(`"Character`", `"extraDamageAfterReload`", `"Flat`", 0.7);`nAddModifier(`"Bullet`", `"maxClip`", `"Flat`", 1);`nAddModifier(`"Grenade`", `"maxClip`", `"Flat`",2);`n} "
$wsAug.Range("J14").Value = "<b>+<color=#c5c5c5ff>50%/70%/<color=#00d100>100%</color></color></b> Damage increase`n<b>+<color=#c5c5c5ff>+1 Bullet, +1 Grenade/+1 Bullet, +2 Grenade/<color=#00d100>+2 Bullet, +2 Grenade</color></color></b> `nThe first shot after a reload does way more damage."
$wsAug.Range("K14").Value = "def OnAttached() {   `nAddModifier(`"Character`", `"extraDamageAfterReload`", `"Flat`", 1);`nAddModifier(`"Bullet`", `"maxClip`", `"Flat`", 2);`nAddModifier(`"Grenade`", `"maxClip`", `"Flat`",2);`n} "

$wsAug.Range("F14:K14").WrapText = $true
$wsAug.Rows.Item(14).RowHeight = 135

# Column widths on the Augments sheet.
$wsAug.Columns.Item(2).ColumnWidth = 13.666666666666666
$wsAug.Columns.Item(5).ColumnWidth = 13.333333333333334
$wsAug.Columns.Item(6).ColumnWidth = 95.16666666666667
$wsAug.Columns.Item(7).ColumnWidth = 45.666666666666664
$wsAug.Columns.Item(8).ColumnWidth = 82.16666666666667
$wsAug.Columns.Item(10).ColumnWidth = 137.33333333333334

# Print setup (picked up a default printer, as happens when Excel
# is used interactively).
$wsAug.PageSetup.PaperSize = 9
$wsAug.PageSetup.Orientation = 1

# Selection / view state left the way it ended up after the edit.
$wsAug.Activate()
$wsAug.Range("J14").Select()

# ----------------------------------------------------------------------
# Sheet "Synergies" (sheet2)
# ----------------------------------------------------------------------

# New row 7: "Dashing" synergy.
$wsSyn.Range("A7").Value = 5
$wsSyn.Range("B7").Value = "Dashing"
$wsSyn.Range("C7").Value = "2,6"
$wsSyn.Range("D7").Value = "icons/Dash.png"
$wsSyn.Range("E7").Value = "Dashing"

# Column widths on the Synergies sheet.
$wsSyn.Columns.Item(4).ColumnWidth = 22.0
$wsSyn.Columns.Item(5).ColumnWidth = 69.66666666666667

$wsSyn.Activate()
$wsSyn.Range("E9").Select()

# Leave the Augments sheet as the active / selected tab.
$wsAug.Activate()
